# Apply the cryptos-list price/volume refresh described in the commit.
# Numeric-looking text values (e.g. "1.00", "5.80") are written via a
# temporary "@" (text) NumberFormat so Excel keeps them as strings instead
# of silently normalizing them into numbers; the format is then reset back
# to the default "Normal" style so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '68.024.95'
$ws.Range("E2").Value = '  +1.36%  '

# Row 3
$ws.Range("D3").Value = '3.269.99'
$ws.Range("E3").Value = '  +0.50%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '587.98'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.81%  '

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '186.86'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +4.80%  '

# Row 7
$ws.Range("E7").Value = '  -0.04%  '

# Row 8
$ws.Range("E8").Value = '  -0.24%  '

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.135'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +3.89%  '

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '6.73'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -0.39%  '

# Row 11
$ws.Range("E11").Value = '  +0.94%  '

# Row 12
$ws.Range("D12").Value = '3.834.38'
$ws.Range("E12").Value = '  +0.45%  '

# Row 13
$ws.Range("E13").Value = '  +0.41%  '

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '28.57'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +1.32%  '

# Row 15
$ws.Range("D15").Value = '68.003.25'
$ws.Range("E15").Value = '  +1.37%  '

# Row 16
$ws.Range("E16").Value = '  +2.45%  '

# Row 17
$ws.Range("D17").Value = '3.267.45'
$ws.Range("E17").Value = '  +0.42%  '

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '5.86'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.10%  '

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '13.61'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +1.45%  '

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '380.77'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +1.62%  '

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '7.73'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +1.04%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.03%  '

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '71.47'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +0.76%  '

# Row 24
$ws.Range("E24").Value = '  +0.32%  '

# Row 25
$ws.Range("E25").Value = '  +1.48%  '

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '9.81'
$c.Style = "Normal"

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.189'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +5.28%  '

# Row 28
$ws.Range("E28").Value = '  +0.02%  '

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '5.80'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +2.65%  '

# Row 30
$ws.Range("E30").Value = '  +0.85%  '

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '22.86'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +1.16%  '

# Row 32
$ws.Range("E32").Value = '  +5.40%  '

# Row 33
$ws.Range("E33").Value = '  +0.01%  '

# Row 34
$ws.Range("E34").Value = '  +0.63%  '

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.54'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +2.28%  '

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '162.66'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -2.44%  '

# Row 37
$ws.Range("E37").Value = '  -0.86%  '

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.838'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -2.32%  '

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '6.80'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +5.03%  '

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '26.60'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -1.45%  '

# Row 41
$ws.Range("E41").Value = '  +4.56%  '

# Row 42
$ws.Range("E42").Value = '  +1.28%  '

# Row 43
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '41.22'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +1.83%  '

# Row 44
$ws.Range("B44").Value = 'Hedera'
$ws.Range("C44").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.0689'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +2.18%  '

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '25.43'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -0.39%  '

# Row 46
$ws.Range("D46").Value = '2.643.70'
$ws.Range("E46").Value = '  -4.48%  '

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '343.13'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -3.18%  '

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.0284'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +1.48%  '

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '32.23'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +4.66%  '

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.997'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +1.04%  '

# Row 51
$ws.Range("E51").Value = '  -0.35%  '
